# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values on the zh-cn and de-de sheets to reflect the latest handback run.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: row 2 (6f39ebe2-... .d47ce87e6def71667d3991212839794704f8b3bc.zh-cn.xlf) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 00:56:49"
$wsZhCn.Range("H2").Value = "2016-03-21 00:57:09"

# --- de-de sheet: row 2 (6f39ebe2-... .d47ce87e6def71667d3991212839794704f8b3bc.de-de.xlf) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 00:56:52"
$wsDeDe.Range("H2").Value = "2016-03-21 00:57:15"
